# Final Demo Check List.docx - add nav grid style notes and bold key terms
#
# 1. Insert two new sub-bullets ("V-for iterations here" and the
#    "@mouseover" note about HomePage.vue / Public speaking) right before
#    the existing "HomePage.vue, PublicSpeaking.vue ..." bullet.
# 2. Bold "style binding" inside the "BlogPage.vue has style binding ..."
#    bullet.

$d = $word.ActiveDocument

$xmlNs = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>{0}</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- locate the anchor paragraph: "HomePage.vue, PublicSpeaking.vue also has similar image switching as assignment 6"
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "HomePage.vue, PublicSpeaking.vue*") {
        $anchorIndex = $i
        break
    }
}

$anchor = $d.Paragraphs.Item($anchorIndex).Range
# Create two blank paragraphs immediately before the anchor paragraph.
$anchor.InsertParagraphBefore()
$anchor.InsertParagraphBefore()

# --- new paragraph 1: "V-for iterations here" (bold "V-for")
$p1Body = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>V-for</w:t></w:r><w:r><w:t xml:space="preserve"> iterations here</w:t></w:r></w:p>'
$target1 = $d.Paragraphs.Item($anchorIndex).Range
$target1.InsertXML(($xmlNs -f $p1Body))

# --- new paragraph 2: "And HomePage.vue (also Public speaking noted below) has @mouseover event bound to image source data attribute in Vue instance (like homework 6) . "
$p2Body = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">And </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>HomePage.vue</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (also Public speaking noted below) has </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>@mouseover</w:t></w:r><w:r><w:t xml:space="preserve"> event bound to image source data attribute in Vue instance (like homework 6</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>) .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$target2 = $d.Paragraphs.Item($anchorIndex + 1).Range
$target2.InsertXML(($xmlNs -f $p2Body))

# --- bold "style binding" inside the BlogPage.vue bullet
$styleIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*has style binding*") {
        $styleIndex = $i
        break
    }
}

$p3 = $d.Paragraphs.Item($styleIndex).Range
$fullText = $p3.Text
$prefixLen = "BlogPage.vue".Length
$subStart = $p3.Start + $prefixLen
$subEnd = $p3.Start + $fullText.Length - 1
$sub = $d.Range($subStart, $subEnd)
$p3Repl = '<w:p><w:r><w:t xml:space="preserve"> has </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>style binding</w:t></w:r><w:r><w:t xml:space="preserve"> (v-bind) for highlighting</w:t></w:r></w:p>'
$sub.InsertXML(($xmlNs -f $p3Repl))

Write-Output "done"
